$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.257729649543762
$ws.Range("B1").Value = 2.390511274337769
$ws.Range("C1").Value = 4.251569747924805
$ws.Range("D1").Value = 2.694173812866211
$ws.Range("E1").Value = 1.354652285575867
